# DXLG yearly financials: a new fiscal-year column was added to the left of
# the existing data (new column D), and the previous D:K columns shifted
# right to E:L. This mirrors selecting column D and using Excel's
# "Insert" (Home > Insert > Insert Sheet Columns) command, then copying the
# number formats from the (old) neighboring column and keying in the new
# year's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column before D; this shifts D:K -> E:L (values,
#    formulas and formatting move with the cells).
$ws.Columns("D:D").Insert()

# 2) The freshly inserted column D has no number formatting yet. Copy the
#    formatting from column E (the column that used to be D) down into D,
#    but only across the row blocks that actually contain data, so we don't
#    manufacture cells on label-only rows (5, 6, 37, 79).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# 3) Fill in the new (most recent) fiscal year's figures in column D.

# -- Income Statement --------------------------------------------------
$ws.Range("D7").Value2 = 43498
$ws.Range("D8").Value2 = 473800
$ws.Range("D9").Value2 = 262500
$ws.Range("D10").Value2 = 211300
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 4600
$ws.Range("D15").Value2 = 28700
$ws.Range("D17").Value2 = "NA"
$ws.Range("D18").Value2 = "NA"
$ws.Range("D20").Value2 = "NA"
$ws.Range("D21").Value2 = "NA"
$ws.Range("D22").Value2 = 3500
$ws.Range("D23").Value2 = -13600
$ws.Range("D24").Value2 = -100
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = "NA"
$ws.Range("D27").Value2 = "NA"
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = "NA"
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = "NA"
$ws.Range("D33").Value2 = "NA"
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = "NA"

# -- Balance Sheet -------------------------------------------------------
$ws.Range("D38").Value2 = 43498
$ws.Range("D41").Value2 = 4900
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 4400
$ws.Range("D44").Value2 = 106800
$ws.Range("D45").Value2 = 11500
$ws.Range("D46").Value2 = 127700
$ws.Range("D47").Value2 = 0
$ws.Range("D48").Value2 = 92500
$ws.Range("D49").Value2 = 1100
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 4700
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 226100
$ws.Range("D57").Value2 = 34400
$ws.Range("D58").Value2 = 41900
$ws.Range("D59").Value2 = 31600
$ws.Range("D60").Value2 = 107900
$ws.Range("D61").Value2 = 14800
$ws.Range("D62").Value2 = 44700
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 167400
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = -153500
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 58600
$ws.Range("D77").Value2 = 0

# -- Cash Flow Statement ---------------------------------------------------
$ws.Range("D80").Value2 = 43498
$ws.Range("D81").Value2 = "NA"
$ws.Range("D83").Value2 = 28700
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 15700
$ws.Range("D91").Value2 = -11800
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -13000
$ws.Range("D96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = -3300
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = -500
